# Applies the Mon Aug 28 03:30:53 UTC 2023 cryptos list refresh (GitHub Actions update).
# Updates Price (col D) and Volume(1h) (col E) figures for the coin rows, and swaps
# the BabyDogeCoin / EnergySwap rows' name, link, price and volume figures (rows 49-50).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold text-formatted numbers/percentages in the
# source sheet (t="inlineStr"). Force the Text number format before writing so the
# COM layer doesn't silently coerce number-looking strings (e.g. "219.37") into
# numeric values, then restore the default (Normal) style once all values are set.
$deRange = $ws.Range("D2:E51")
$deRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.159.05"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "1.658.80"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").Value = "219.37"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").Value = "0.5227"
$ws.Range("E6").Value = "  -1.59%  "
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "0.06299"
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("D10").Value = "20.64"
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("D11").Value = "0.07818"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "4.495"
$ws.Range("E12").Value = "  -1.51%  "
$ws.Range("D13").Value = "1.656.48"
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("D14").Value = "1.886.90"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "0.5550"
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").Value = "0.0₅8013"
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("D17").Value = "65.12"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").Value = "26.189.76"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").Value = "4.639"
$ws.Range("D21").Value = "196.12"
$ws.Range("E21").Value = "  +1.76%  "
$ws.Range("D22").Value = "10.11"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").Value = "5.962"
$ws.Range("E23").Value = "  -0.97%  "
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").Value = "146.35"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("D26").Value = "0.1202"
$ws.Range("E26").Value = "  -1.48%  "
$ws.Range("D27").Value = "7.143"
$ws.Range("E27").Value = "  -0.53%  "
$ws.Range("D28").Value = "16.00"
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("D29").Value = "1.493"
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("D30").Value = "0.05742"
$ws.Range("E30").Value = "  -2.28%  "
$ws.Range("D31").Value = "1.276"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("D32").Value = "3.494"
$ws.Range("D33").Value = "3.372"
$ws.Range("E33").Value = "  +3.30%  "
$ws.Range("D34").Value = "1.586"
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("D35").Value = "0.9556"
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("E36").Value = "  -0.62%  "
$ws.Range("D37").Value = "2.420"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").Value = "0.5722"
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("D39").Value = "0.01596"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").Value = "5.962"
$ws.Range("E40").Value = "  +2.27%  "
$ws.Range("D41").Value = "1.064.96"
$ws.Range("E41").Value = "  +2.17%  "
$ws.Range("D42").Value = "0.8493"
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").Value = "103.94"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "1.797.38"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("D46").Value = "58.05"
$ws.Range("E46").Value = "  +1.25%  "
$ws.Range("D47").Value = "1.005"
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "8.041"
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₈102"
$ws.Range("E50").Value = "  -3.10%  "
$ws.Range("D51").Value = "0.05202"
$ws.Range("E51").Value = "  +0.73%  "

$deRange.Style = "Normal"

Write-Output "Applied 94 cell updates to cryptos sheet"
